# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sat Jun 24 23:25:44 UTC 2023 with GitHub Actions".
# Column D = Price, column E = Volume(1h) percentage change. Both columns
# hold plain text in the source sheet (e.g. "30.502.57", "  -0.56%  "), so
# each cell is forced to Text format before the write (and reset back to
# the default "Normal" style afterwards) to stop Excel from reinterpreting
# the numeric-looking strings as actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the sheet row plus the new Price (column D) and
# Volume(1h) (column E) text. $null means that column is unchanged.
$updates = @(
    @{ Row = 2; D = "30.534.15"; E = "  -0.32%  " }
    @{ Row = 3; D = "1.875.01"; E = "  -0.97%  " }
    @{ Row = 4; D = "0.9991"; E = "  -0.20%  " }
    @{ Row = 5; D = "236.44"; E = "  -3.30%  " }
    @{ Row = 6; D = "0.9992"; E = "  -0.14%  " }
    @{ Row = 7; D = "0.4871"; E = "  -2.11%  " }
    @{ Row = 8; D = "0.2900"; E = "  -2.04%  " }
    @{ Row = 9; D = "0.06669"; E = "  -2.27%  " }
    @{ Row = 10; D = "1.871.43"; E = "  -1.15%  " }
    @{ Row = 11; D = "16.60"; E = "  -2.72%  " }
    @{ Row = 12; D = "0.07225"; E = "  -1.52%  " }
    @{ Row = 13; D = "89.15"; E = "  -2.52%  " }
    @{ Row = 14; D = "5.010"; E = "  -1.75%  " }
    @{ Row = 15; D = "0.6529"; E = "  -3.17%  " }
    @{ Row = 16; D = "30.466.06"; E = "  -0.56%  " }
    @{ Row = 17; D = "0.000007847"; E = "  -0.91%  " }
    @{ Row = 18; D = "0.9991"; E = "  -0.13%  " }
    @{ Row = 19; D = "13.01"; E = "  -2.05%  " }
    @{ Row = 20; D = "2.112.94"; E = "  -0.95%  " }
    @{ Row = 21; D = "0.9995"; E = "  -0.16%  " }
    @{ Row = 22; D = "213.75"; E = "  +19.72%  " }
    @{ Row = 23; D = "4.726"; E = "  -3.00%  " }
    @{ Row = 24; D = "6.131"; E = "  +0.95%  " }
    @{ Row = 25; D = "9.372"; E = "  +0.78%  " }
    @{ Row = 26; D = "156.13"; E = "  +1.15%  " }
    @{ Row = 27; D = "18.79"; E = "  +0.04%  " }
    @{ Row = 28; D = "1.828"; E = "  -5.13%  " }
    @{ Row = 29; D = "1.406"; E = "  +1.38%  " }
    @{ Row = 30; D = "4.260"; E = "  -1.81%  " }
    @{ Row = 31; D = "0.09026"; E = "  +0.87%  " }
    @{ Row = 32; D = "3.925"; E = "  -2.76%  " }
    @{ Row = 33; D = "0.05100"; E = "  -2.05%  " }
    @{ Row = 34; D = "0.7234"; E = "  -2.32%  " }
    @{ Row = 35; D = "1.078"; E = "  -5.10%  " }
    @{ Row = 36; D = "2.686"; E = "  +0.49%  " }
    @{ Row = 37; D = "0.01810"; E = "  -3.47%  " }
    @{ Row = 39; D = "0.9185"; E = $null }
    @{ Row = 40; D = "2.042"; E = "  -6.02%  " }
    @{ Row = 41; D = "0.4388"; E = "  +0.64%  " }
    @{ Row = 42; D = "104.45"; E = "  -1.51%  " }
    @{ Row = 43; D = "5.741"; E = "  -1.16%  " }
    @{ Row = 44; D = "0.9945"; E = "  -0.63%  " }
    @{ Row = 45; D = "0.1326"; E = "  -2.10%  " }
    @{ Row = 46; D = "7.327"; E = "  -4.37%  " }
    @{ Row = 47; D = "0.4011"; E = "  +3.03%  " }
    @{ Row = 48; D = "0.05822"; E = "  -0.44%  " }
    @{ Row = 49; D = "8.658"; E = "  +2.03%  " }
    @{ Row = 50; D = "1.403"; E = "  +1.56%  " }
    @{ Row = 51; D = "33.12"; E = "  -0.97%  " }
)

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

foreach ($u in $updates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    if ($null -ne $u.E) {
        Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
    }
}
